$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 112154283
$ws.Range("AO9").Value = "på blmr av åkervädd"
$ws.Range("B9").Value = 44331
$ws.Range("E9").Value = 201164
$ws.Range("F9").Value = "Sexfläckig bastardsvärmare"
$ws.Range("G9").Value = "Zygaena filipendulae"
$ws.Range("I9").Value = "'1"
$ws.Range("M9").Value = "födosökande"
$ws.Range("Q9").Value = 442664.1890363992
$ws.Range("R9").Value = 6204260.315617888

# Row 10
$ws.Range("AO10").ClearContents()
$ws.Range("A10").Value = 112154281
$ws.Range("B10").Value = 42578
$ws.Range("E10").Value = 100770
$ws.Range("F10").Value = "Mindre blåvinge"
$ws.Range("G10").Value = "Cupido minimus"
$ws.Range("H10").Value = "(Fuessly, 1775)"
$ws.Range("I10").Value = "'1"
$ws.Range("M10").Value = "friflygande"
$ws.Range("Q10").Value = 442664.1890363992
$ws.Range("R10").Value = 6204260.315617888

# Row 11
$ws.Range("A11").Value = 112154275
$ws.Range("AO11").Value = "på blmr av åkervädd m fl"
$ws.Range("B11").Value = 44322
$ws.Range("E11").Value = 102366
$ws.Range("F11").Value = "Ängsmetallvinge"
$ws.Range("G11").Value = "Adscita statices"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("I11").Value = "'5"
$ws.Range("M11").Value = "vilande"
$ws.Range("Q11").Value = 442616.138687243
$ws.Range("R11").Value = 6204441.08982533

# Row 12
$ws.Range("A12").Value = 112154282
$ws.Range("B12").Value = 44322
$ws.Range("E12").Value = 102366
$ws.Range("F12").Value = "Ängsmetallvinge"
$ws.Range("G12").Value = "Adscita statices"
$ws.Range("I12").Value = "'4"
$ws.Range("M12").Value = "vilande"

# Row 14
$ws.Range("A14").Value = 112145591
$ws.Range("AA14").Value = "'2013-06-12"
$ws.Range("AC14").Value = "lufthåvning"
$ws.Range("AF14").Value = "'"
$ws.Range("AI14").Value = "i igenväxande hed"
$ws.Range("AO14").Value = "på tjärblomster"
$ws.Range("AQ14").Value = "Nils Otto Nilsson"
$ws.Range("AR14").Value = "NON 04542"
$ws.Range("AX14").Value = "Nils Otto Nilsson"
$ws.Range("B14").Value = 42546
$ws.Range("E14").Value = 102923
$ws.Range("F14").Value = "Violettkantad guldvinge"
$ws.Range("G14").Value = "Lycaena hippothoe"
$ws.Range("H14").Value = "(Linnaeus, 1760)"
$ws.Range("K14").Value = "imago/adult"
$ws.Range("L14").Value = "hane"
$ws.Range("M14").Value = "födosökande"
$ws.Range("N14").Value = "'"
$ws.Range("Q14").Value = 442542.9522655545
$ws.Range("R14").Value = 6204459.965453062
$ws.Range("Y14").Value = "'2013-06-12"

# Row 15
$ws.Range("A15").Value = 112154276
$ws.Range("B15").Value = 39449
$ws.Range("E15").Value = 102471
$ws.Range("F15").Value = "Åkerväddsantennmal"
$ws.Range("G15").Value = "Nemophora metallica"
$ws.Range("H15").Value = "(Poda, 1761)"
$ws.Range("I15").Value = "'1"
$ws.Range("Q15").Value = 442616.138687243
$ws.Range("R15").Value = 6204441.08982533

# Row 17
$ws.Range("A17").Value = 112154273
$ws.Range("AA17").Value = "'2013-07-09"
$ws.Range("AB17").Value = "'00:00"
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AI17").Value = "på igenväxande sandhed"
$ws.Range("AO17").Value = "på blmr av åkervädd m fl"
$ws.Range("AT17").Value = "'"
$ws.Range("AW17").Value = "Nils Otto Nilsson"
$ws.Range("AX17").Value = "Nils Otto Nilsson, Mats Karlsson"
$ws.Range("AY17").Value = "Krst NV-program 2013"
$ws.Range("B17").Value = 44322
$ws.Range("C17").Value = "Ovaliderad"
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 102366
$ws.Range("F17").Value = "Ängsmetallvinge"
$ws.Range("G17").Value = "Adscita statices"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("I17").Value = "'4"
$ws.Range("J17").Value = "ex."
$ws.Range("M17").Value = "vilande"
$ws.Range("P17").Value = "Vä, delomr 22, 700 m NNO Sigridslund, Sk"
$ws.Range("Q17").Value = 442541.7120545401
$ws.Range("R17").Value = 6204451.031370129
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = "Skåne"
$ws.Range("U17").Value = "Kristianstad"
$ws.Range("V17").Value = "Skåne"
$ws.Range("W17").Value = "Vä"
$ws.Range("Y17").Value = "'2013-07-09"
$ws.Range("Z17").Value = "'00:00"
